$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newWs.Name = "Select Values"

$ws2 = $wb.Worksheets.Item("ETA calculator")
$ws2.Activate()
